$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("B2").Value = "CS303"
$wsA.Range("C2").Value = "CS309 (Tutorial)"
$wsA.Range("D2").Value = "Free"
$wsA.Range("B3").Value = "Free"
$wsA.Range("F3").Value = "CS309"
$wsA.Range("C5").Value = "CS304"
$wsA.Range("D5").Value = "CS303"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "CS304 (Tutorial)"
$wsA.Range("B6").Value = "CS304"
$wsA.Range("C6").Value = "CS303"
$wsA.Range("D6").Value = "CS303 (Tutorial)"
$wsA.Range("E6").Value = "CS461 (Elective)"
$wsA.Range("F6").Value = "Free"
$wsA.Range("C7").Value = "Free"
$wsA.Range("D7").Value = "CS309"
$wsA.Range("E7").Value = "CS309"
$wsA.Range("F7").Value = "Free"

$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "Free"
$wsB.Range("D2").Value = "CS309 (Tutorial)"
$wsB.Range("E2").Value = "CS304 (Tutorial)"
$wsB.Range("D3").Value = "CS309"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "CS309"
$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "CS304"
$wsB.Range("E5").Value = "CS303"
$wsB.Range("F5").Value = "CS304"
$wsB.Range("B6").Value = "CS309"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "CS461 (Elective)"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "CS303"
$wsB.Range("E7").Value = "CS304"
$wsB.Range("F7").Value = "CS303 (Tutorial)"
